$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7: update description text, keep test-case list but extend it ---
$ws.Range("C7").Value = "Check the Home page Functionality"
$ws.Range("D7").Value = "1. All Nav Menu links navigate to the proper area of the page`n2. All Header Links are working`n3. All Footer Links are working`n4. All Main Page Links are working`n5. All Images are working"

# --- Row 8 ---
$ws.Range("C8").Value = "Check the Fireplace Sets page Functionality "
$ws.Range("D8").Value = "1. All Nav Menu links navigate to the proper area of the page`n2. All Header Links are working`n3. All Footer Links are working`n4. All Main Page Links are working`n5. All Images are working"
$ws.Range("D8").WrapText = $true

# --- Row 9 ---
$ws.Range("C9").Value = "Check the Tank Systems Page"
$ws.Range("D9").Value = "1. All Nav Menu links navigate to the proper area of the page`n2. All Header Links are working`n3. All Footer Links are working`n4. All Main Page Links are working`n5. All Images are working"
$ws.Range("D9").WrapText = $true

# --- Row 10 ---
$ws.Range("C10").Value = "Check the Instructions Functionality"
$ws.Range("D10").Value = "1. All Nav Menu links navigate to the proper area of the page`n2. All Header Links are working`n3. All Footer Links are working`n4. All Main Page Links are working`n5. All Images are working"
$ws.Range("D10").WrapText = $true

# --- Row 11 ---
$ws.Range("C11").Value = "Check the Q&A Functionality"
$ws.Range("D11").Value = "1. All Nav Menu links navigate to the proper area of the page`n2. All Header Links are working`n3. All Footer Links are working`n4. All Main Page Links are working`n5. All Images are working"
$ws.Range("D11").WrapText = $true

# --- Row 12 ---
$ws.Range("C12").Value = "Check Specials Page Functionality"
$ws.Range("D12").Value = "1. Check that the shopify page loads`n2. All Nav Menu links navigate to the proper page"
$ws.Range("D12").WrapText = $true

# --- Row 13 ---
$ws.Range("C13").Value = "Check the Website Funcitonality on Mobile Phone"
$ws.Range("D13").Value = "1. Dropdown Nav Menu works properly`n2. Header Nav Menu links are working`n3. Footer menu links are working`n4. Speed Links are working"
$ws.Range("D13").WrapText = $true

# --- Row 14 & 15: clear the now-unused test scenarios ---
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()

# --- Row heights ---
$ws.Rows.Item(7).RowHeight = 75
$ws.Rows.Item(8).RowHeight = 75
$ws.Rows.Item(9).RowHeight = 75
$ws.Rows.Item(10).RowHeight = 75
$ws.Rows.Item(11).RowHeight = 75
$ws.Rows.Item(12).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 60

# --- Column C width (closest representable value to the target 45.5703125) ---
$ws.Columns.Item(3).ColumnWidth = 44.65

# --- Selection ---
$ws.Range("C3").Select() | Out-Null
